# Daily attendance processing - reorder the "Recorded By" contributor
# lists in column G so that they are listed in reverse order.
#
# Only touches cells whose value is a comma-separated list made up
# exclusively of the known "recorder" tokens (System / system,
# backup@backdoor.com, dnasr281@gmail.com). Cells containing other
# recorders (e.g. admin@admin.com) or a single value are left as-is,
# matching the source data exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$allowedTokens = @("System", "system", "backup@backdoor.com", "dnasr281@gmail.com")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value) { continue }
    if (-not ($value -is [string])) { continue }
    if ($value.IndexOf(",") -lt 0) { continue }

    $parts = $value -split ", "

    $allKnown = $true
    foreach ($part in $parts) {
        if ($allowedTokens -notcontains $part) {
            $allKnown = $false
            break
        }
    }

    if (-not $allKnown) { continue }

    $reversedParts = @()
    for ($i = $parts.Count - 1; $i -ge 0; $i--) {
        $reversedParts += $parts[$i]
    }
    $newValue = [string]::Join(", ", $reversedParts)

    if ($newValue -ne $value) {
        $cell.Value2 = $newValue
    }
}
